#
# Scheduled market-data refresh: update current Market Board price snapshots
# (columns H/I/J/K/L) and the derived Leve-profit columns (M/N) on each job sheet.
#
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(51, 8).Value = 15631861
$ws.Cells.Item(51, 10).Value = 8482.833000000001
$ws.Cells.Item(51, 12).Value = 8482.833000000001
$ws.Cells.Item(51, 14).Value = -9450.833000000001

$ws.Cells.Item(116, 8).Value = 13252.125
$ws.Cells.Item(116, 9).Value = 15025.134
$ws.Cells.Item(116, 11).Value = 15025.134
$ws.Cells.Item(116, 13).Value = -11583.134


$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 34304.883
$ws.Cells.Item(32, 9).Value = 39495.414
$ws.Cells.Item(32, 10).Value = 4199.8
$ws.Cells.Item(32, 11).Value = 39495.414
$ws.Cells.Item(32, 12).Value = 4199.8
$ws.Cells.Item(32, 13).Value = -39208.414
$ws.Cells.Item(32, 14).Value = -4773.8

$ws.Cells.Item(55, 8).Value = 5000
$ws.Cells.Item(55, 9).Value = 5000
$ws.Cells.Item(55, 11).Value = 5000
$ws.Cells.Item(55, 13).Value = -4685

$ws.Cells.Item(80, 8).Value = 32500.25

$ws.Cells.Item(83, 8).Value = 32500.25


$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 870.9286
$ws.Cells.Item(80, 9).Value = 851.6
$ws.Cells.Item(80, 11).Value = 851.6
$ws.Cells.Item(80, 13).Value = 146.4

$ws.Cells.Item(83, 8).Value = 870.9286
$ws.Cells.Item(83, 9).Value = 851.6
$ws.Cells.Item(83, 11).Value = 4258
$ws.Cells.Item(83, 13).Value = 734

$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 14).Value = ""

$ws.Cells.Item(91, 8).Value = 0
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 14).Value = ""

$ws.Cells.Item(99, 8).Value = 1521.9286
$ws.Cells.Item(99, 9).Value = 1358.9166
$ws.Cells.Item(99, 11).Value = 1358.9166
$ws.Cells.Item(99, 13).Value = 139.0834

$ws.Cells.Item(107, 8).Value = 37550.93
$ws.Cells.Item(107, 9).Value = 45881.91
$ws.Cells.Item(107, 11).Value = 45881.91
$ws.Cells.Item(107, 13).Value = -43961.91


$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(10, 8).Value = 419
$ws.Cells.Item(10, 9).Value = 419
$ws.Cells.Item(10, 11).Value = 419
$ws.Cells.Item(10, 13).Value = -280

$ws.Cells.Item(47, 8).Value = 1000
$ws.Cells.Item(47, 9).Value = 1000
$ws.Cells.Item(47, 11).Value = 1000
$ws.Cells.Item(47, 13).Value = -434

$ws.Cells.Item(99, 8).Value = 8711.294
$ws.Cells.Item(99, 9).Value = 5109.3335
$ws.Cells.Item(99, 10).Value = 12763.5
$ws.Cells.Item(99, 11).Value = 5109.3335
$ws.Cells.Item(99, 12).Value = 12763.5
$ws.Cells.Item(99, 13).Value = -3611.3335
$ws.Cells.Item(99, 14).Value = -15759.5

$ws.Cells.Item(126, 8).Value = 8711.294
$ws.Cells.Item(126, 9).Value = 5109.3335
$ws.Cells.Item(126, 10).Value = 12763.5
$ws.Cells.Item(126, 11).Value = 15328.0005
$ws.Cells.Item(126, 12).Value = 38290.5
$ws.Cells.Item(126, 13).Value = -12858.0005
$ws.Cells.Item(126, 14).Value = -43230.5

$ws.Cells.Item(132, 8).Value = 53871.633
$ws.Cells.Item(132, 9).Value = 59856.53
$ws.Cells.Item(132, 11).Value = 179569.59
$ws.Cells.Item(132, 13).Value = -177039.59

$ws.Cells.Item(134, 8).Value = 2069.5293
$ws.Cells.Item(134, 9).Value = 1598.6666
$ws.Cells.Item(134, 11).Value = 4795.9998
$ws.Cells.Item(134, 13).Value = -2260.9998


$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(102, 8).Value = 0
$ws.Cells.Item(102, 9).Value = 0
$ws.Cells.Item(102, 11).Value = 0

$ws.Cells.Item(103, 8).Value = 1201.5714
$ws.Cells.Item(103, 9).Value = 700
$ws.Cells.Item(103, 10).Value = 1402.2
$ws.Cells.Item(103, 11).Value = 2100
$ws.Cells.Item(103, 12).Value = 4206.6
$ws.Cells.Item(103, 13).Value = -1221
$ws.Cells.Item(103, 14).Value = -5964.6

$ws.Cells.Item(104, 8).Value = 0
$ws.Cells.Item(104, 9).Value = 0
$ws.Cells.Item(104, 11).Value = 0

$ws.Cells.Item(105, 8).Value = 7199
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 12).Value = 0


$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1587.2916
$ws.Cells.Item(97, 9).Value = 1022.3333
$ws.Cells.Item(97, 10).Value = 2528.889
$ws.Cells.Item(97, 11).Value = 1022.3333
$ws.Cells.Item(97, 12).Value = 2528.889
$ws.Cells.Item(97, 13).Value = -526.3333
$ws.Cells.Item(97, 14).Value = -3520.889

$ws.Cells.Item(122, 8).Value = 2755.6365
$ws.Cells.Item(122, 9).Value = 2312.9
$ws.Cells.Item(122, 10).Value = 3124.5833
$ws.Cells.Item(122, 11).Value = 6938.700000000001
$ws.Cells.Item(122, 12).Value = 9373.749899999999
$ws.Cells.Item(122, 13).Value = -4488.700000000001
$ws.Cells.Item(122, 14).Value = -14273.7499

$ws.Cells.Item(126, 8).Value = 2911.0454
$ws.Cells.Item(126, 9).Value = 2003.6154
$ws.Cells.Item(126, 11).Value = 6010.8462
$ws.Cells.Item(126, 13).Value = -3540.8462


$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3664.5
$ws.Cells.Item(7, 9).Value = 2996.4
$ws.Cells.Item(7, 11).Value = 2996.4
$ws.Cells.Item(7, 13).Value = -2884.4

$ws.Cells.Item(22, 8).Value = 2140
$ws.Cells.Item(22, 10).Value = 2175
$ws.Cells.Item(22, 12).Value = 2175
$ws.Cells.Item(22, 14).Value = -2765

$ws.Cells.Item(27, 8).Value = 2140
$ws.Cells.Item(27, 10).Value = 2175
$ws.Cells.Item(27, 12).Value = 2175
$ws.Cells.Item(27, 14).Value = -2389

$ws.Cells.Item(40, 8).Value = 2374.3333
$ws.Cells.Item(40, 9).Value = 2226.5454
$ws.Cells.Item(40, 11).Value = 2226.5454
$ws.Cells.Item(40, 13).Value = -2090.5454

$ws.Cells.Item(46, 8).Value = 4691.3
$ws.Cells.Item(46, 9).Value = 749.6667
$ws.Cells.Item(46, 11).Value = 749.6667
$ws.Cells.Item(46, 13).Value = -561.6667

$ws.Cells.Item(55, 8).Value = 732.6799999999999
$ws.Cells.Item(55, 9).Value = 246.91667
$ws.Cells.Item(55, 10).Value = 1181.0769
$ws.Cells.Item(55, 11).Value = 246.91667
$ws.Cells.Item(55, 12).Value = 1181.0769
$ws.Cells.Item(55, 13).Value = -73.91667000000001
$ws.Cells.Item(55, 14).Value = -1527.0769

$ws.Cells.Item(93, 8).Value = 1891.6154
$ws.Cells.Item(93, 9).Value = 1849.1111
$ws.Cells.Item(93, 10).Value = 1987.25
$ws.Cells.Item(93, 11).Value = 1849.1111
$ws.Cells.Item(93, 12).Value = 1987.25
$ws.Cells.Item(93, 13).Value = -601.1111000000001
$ws.Cells.Item(93, 14).Value = -4483.25

$ws.Cells.Item(100, 8).Value = 3333
$ws.Cells.Item(100, 9).Value = 2800
$ws.Cells.Item(100, 11).Value = 2800
$ws.Cells.Item(100, 13).Value = -2259

$ws.Cells.Item(126, 8).Value = 3664.5
$ws.Cells.Item(126, 9).Value = 2996.4
$ws.Cells.Item(126, 11).Value = 8989.200000000001
$ws.Cells.Item(126, 13).Value = -6519.200000000001

$ws.Cells.Item(136, 8).Value = 3124.7
$ws.Cells.Item(136, 9).Value = 3124.7
$ws.Cells.Item(136, 11).Value = 9374.099999999999
$ws.Cells.Item(136, 13).Value = -6824.099999999999


$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 281510.84
$ws.Cells.Item(126, 9).Value = 2759.1
$ws.Cells.Item(126, 11).Value = 8277.299999999999
$ws.Cells.Item(126, 13).Value = -5807.299999999999
